$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 103
$ws.Range("B2").Value = "Eduardo"
$ws.Range("C2").Value = "Educado"

$ws.Range("A3").Value = 105
$ws.Range("B3").Value = "Lola"
$ws.Range("C3").Value = "Sol"
